$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 65 (F, G, I change) ---
$ws.Range("F65").Value = 'Dagang Bahan Bangunan'
$ws.Range("G65").Value = 'KOTA JUANG'
$ws.Range("I65").Value = '09 November 2017'

# --- Insert 11 new rows at 66..76 (pushes old rows 66+ down by 11) ---
$ws.Range("A66:A76").EntireRow.Insert()

# --- Apply formatting to the newly inserted rows to match the report style ---
$dataRange = $ws.Range("A66:I76")
$dataRange.RowHeight = 30
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true
$dataRange.Borders.LineStyle = 1

$markRange = $ws.Range("J66:L76")
$markRange.HorizontalAlignment = -4108
$markRange.VerticalAlignment = -4108
$markRange.WrapText = $true
$markRange.Borders.LineStyle = 1

# --- Fill in the data for the new rows ---
# Row 66
$ws.Range("A66").Value = 58
$ws.Range("B66").Value = '510.1/0059/KPPTSP/2015'
$ws.Range("C66").Value = 'JAFAR ISHAK'
$ws.Range("D66").Value = ' "METRO TANI"'
$ws.Range("E66").Value = 'Jl. Jangka Dsn. Suka Mulia Gp. Geundot'
$ws.Range("F66").Value = 'Dagang Pupuk dan Pestisida'
$ws.Range("G66").Value = 'JANGKA'
$ws.Range("H66").Value = '23 Januari 2015'
$ws.Range("I66").Value = '11 Februari 2018'
$ws.Range("K66").Value = '√'

# Row 67
$ws.Range("A67").Value = 59
$ws.Range("B67").Value = '510.1/0061/KPPTSP/2015'
$ws.Range("C67").Value = 'MARTUNIS'
$ws.Range("D67").Value = ' "UD. JEUMPA PUTEH"'
$ws.Range("E67").Value = 'Jl. Medan - B. Aceh / Limeng Madoe Gp. Cot Bada'
$ws.Range("F67").Value = 'Dagang Pupuk dan Pestisida'
$ws.Range("G67").Value = 'JEUMPA'
$ws.Range("H67").Value = '23 Januari 2015'
$ws.Range("I67").Value = '08 Februari 2018'
$ws.Range("K67").Value = '√'

# Row 68
$ws.Range("A68").Value = 60
$ws.Range("B68").Value = '510.1/0062/KPPTSP/2015'
$ws.Range("C68").Value = 'HAMZAH'
$ws.Range("D68").Value = ' "UD. DARA DITA"'
$ws.Range("E68").Value = ' Gp. Jarommah Baroh'
$ws.Range("F68").Value = 'Dagang Bahan Kebutuhan Pokok'
$ws.Range("G68").Value = 'KUTA BLANG'
$ws.Range("H68").Value = '26 Januari 2015'
$ws.Range("I68").Value = '25 Januari 2018'
$ws.Range("J68").Value = '√'

# Row 69
$ws.Range("A69").Value = 61
$ws.Range("B69").Value = '510.1/0063/KPPTSP/2015'
$ws.Range("C69").Value = 'MUSTAHAR'
$ws.Range("D69").Value = ' "CV. COBRA JAYA"'
$ws.Range("E69").Value = 'Jl. B. Aceh - Medan Gp. Matang Sagoe'
$ws.Range("F69").Value = 'Kontraktor - Leveransir'
$ws.Range("G69").Value = 'PEUSANGAN'
$ws.Range("H69").Value = '26 Januari 2015'
$ws.Range("I69").Value = '01 Februari 2018'
$ws.Range("K69").Value = '√'

# Row 70
$ws.Range("A70").Value = 62
$ws.Range("B70").Value = '510.1/0064/KPPTSP/2015'
$ws.Range("C70").Value = 'MUNAWIR'
$ws.Range("D70").Value = ' "CV. JULI RAYA"'
$ws.Range("E70").Value = 'Jl. Takengon Km. 6 Gp. Blang Keutumba'
$ws.Range("F70").Value = 'Kontraktor - Leveransir'
$ws.Range("G70").Value = 'JULI'
$ws.Range("H70").Value = '26 Januari 2015'
$ws.Range("I70").Value = '25 Januari 2018'
$ws.Range("J70").Value = '√'

# Row 71
$ws.Range("A71").Value = 63
$ws.Range("B71").Value = '510.1/0067/KPPTSP/2015'
$ws.Range("C71").Value = 'MUTTAQIN, S.Farm, Apt'
$ws.Range("D71").Value = ' "APOTIK JAKARTA"'
$ws.Range("E71").Value = 'Jl. Medan - B. Aceh Gp. Tingkeum Manyang'
$ws.Range("F71").Value = 'Apotik'
$ws.Range("G71").Value = 'KUTA BLANG'
$ws.Range("H71").Value = '27 Januari 2015'
$ws.Range("I71").Value = '26 Januari 2018'
$ws.Range("J71").Value = '√'

# Row 72
$ws.Range("A72").Value = 64
$ws.Range("B72").Value = '510.1/0065/KPPTSP/2015'
$ws.Range("C72").Value = 'MUNAWIR'
$ws.Range("D72").Value = ' "KAYLA BUTIK ONLINE"'
$ws.Range("E72").Value = 'Jl. Medan - B. Aceh Gp. Matang Glp. Dua Mns. Dayah'
$ws.Range("F72").Value = 'Jualan Pakaian Jadi'
$ws.Range("G72").Value = 'PEUSANGAN'
$ws.Range("H72").Value = '26 Januari 2015'
$ws.Range("I72").Value = '25 Januari 2018'
$ws.Range("J72").Value = '√'

# Row 73
$ws.Range("A73").Value = 65
$ws.Range("B73").Value = '510.1/0066/KPPTSP/2015'
$ws.Range("C73").Value = 'SYAHRIL M. DAUD'
$ws.Range("D73").Value = ' "PT. KONSTRUKSI BUMI NUSANTARA"'
$ws.Range("E73").Value = 'Jln. B. Aceh - Medan No. 20-21  Gp. Cot Gapu'
$ws.Range("F73").Value = 'Kontraktor - Leveransir'
$ws.Range("G73").Value = 'KOTA JUANG'
$ws.Range("H73").Value = '27 Januari 2015'
$ws.Range("I73").Value = '15 November 2017'
$ws.Range("K73").Value = '√'

# Row 74
$ws.Range("A74").Value = 66
$ws.Range("B74").Value = '510.1/0068/KPPTSP/2015'
$ws.Range("C74").Value = 'ZULKIFLI IBRAHIM'
$ws.Range("D74").Value = ' "CV. PRIMA KECANA"'
$ws.Range("E74").Value = 'Jl. B. Aceh - Medan Gp. Bandar Bireuen'
$ws.Range("F74").Value = 'Kontraktor - Leveransir'
$ws.Range("G74").Value = 'KOTA JUANG'
$ws.Range("H74").Value = '27 Januari 2015'
$ws.Range("I74").Value = '28 September 2017'
$ws.Range("K74").Value = '√'

# Row 75
$ws.Range("A75").Value = 67
$ws.Range("B75").Value = '510.1/0069/KPPTSP/2015'
$ws.Range("C75").Value = 'NUZUL IHSAN'
$ws.Range("D75").Value = ' "UD. PERMATA GAS"'
$ws.Range("E75").Value = 'Jl. Buket Teukuh Gp. Bukit Teukueh'
$ws.Range("F75").Value = 'Pangkalan LPG 3 Kg (Subsidi)'
$ws.Range("G75").Value = 'KOTA JUANG'
$ws.Range("H75").Value = '28 Januari 2015'
$ws.Range("I75").Value = '27 Januari 2018'
$ws.Range("J75").Value = '√'

# Row 76
$ws.Range("A76").Value = 68
$ws.Range("B76").Value = '510.1/0070/KPPTSP/2015'
$ws.Range("C76").Value = 'HAMZATUL IQBAL'
$ws.Range("D76").Value = ' "TUNAS BARU SERVICE"'
$ws.Range("E76").Value = 'Jl. Sinar Peusangan No. 23  Gp. Keude Matang Glumpang Dua'
$ws.Range("F76").Value = 'Dagang Spare Part Honda dan Service'
$ws.Range("G76").Value = 'PEUSANGAN'
$ws.Range("H76").Value = '28 Januari 2015'
$ws.Range("I76").Value = '27 Januari 2018'
$ws.Range("J76").Value = '√'

$ws.Range("L76").Select()
